$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Parameter block (A1:B4) ---
$ws.Range("A1").Value = "min"
$ws.Range("B1").Value = 0

$ws.Range("A2").Value = "max"
$ws.Range("B2").Value = 100

$ws.Range("A4").Value = "type of scale"
$ws.Range("B4").Value = "linear"

$ws.Range("A3").Value = "number of bins"
$ws.Range("B3").Value = 10

# --- Bin counter row (row 8) ---
$ws.Range("E8").Value = "bins"
$ws.Range("F8").Value = 1
$ws.Range("G8").Formula = "=IF(F`$8<`$B`$3,F`$8+1,"""")"
$ws.Range("H8:Y8").Formula = "=IF(G`$8<`$B`$3,G`$8+1,"""")"

# --- Values (low) row (row 9) ---
$ws.Range("E9").Value = "Values (low)"
$ws.Range("F9").Formula = "=`$B`$1"

# --- Scale list (rows 19-21) ---
$ws.Range("A19").Value = "scales"
$ws.Range("A20").Value = "linear"
$ws.Range("A21").Value = "log"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 23.88671875
$ws.Columns.Item(5).ColumnWidth = 10.88671875
$ws.Columns.Item(6).ColumnWidth = 10.109375

# --- Data validation on B4: list restricted to A20:A21 ---
$validation = $ws.Range("B4").Validation
$validation.Delete()
$validation.Add(3, 1, 1, "=`$A`$20:`$A`$21")
$validation.ShowInput = $true
$validation.ShowError = $true
$validation.InputMessage = "choose a linear or a log scale"

# --- Selection matches diff (active cell F9) ---
$ws.Range("F9").Select()

$wb.Save()
